$d = $word.ActiveDocument

# Locate, by content, the paragraph "LOB1018: Física I (Requisito fraco)"
# and the paragraph starting with the "© 2020 ... Contact: luizeleno@usp.br"
# copyright notice. Everything strictly between them (the blank paragraph,
# the "Ver no Jupiter..." paragraph, and that copyright paragraph itself)
# must be removed, while the "LOB1018..." paragraph is kept untouched.

$startIdx = -1
$endIdx = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "*LOB1018: Física I*") {
        $startIdx = $i
    }
    if ($txt -like "*Contact: luizeleno@usp.br*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1 -and $endIdx -gt $startIdx) {
    $rangeStart = $d.Paragraphs.Item($startIdx + 1).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIdx).Range.End
    $victim = $d.Range($rangeStart, $rangeEnd)
    $victim.Delete()
}
